$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing it to stay text,
# since some "Price" values (e.g. "569.57") would otherwise be
# auto-converted into numbers by Excel on assignment.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "64.452.31"
$ws.Cells.Item(2, 5).Value = "  -2.15%  "

Set-TextValue 3 4 "3.413.26"
$ws.Cells.Item(3, 5).Value = "  -2.56%  "

$ws.Cells.Item(4, 5).Value = "  +0.04%  "

Set-TextValue 5 4 "569.57"
$ws.Cells.Item(5, 5).Value = "  -2.04%  "

Set-TextValue 6 4 "157.25"
$ws.Cells.Item(6, 5).Value = "  -2.88%  "

$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue 7 4 "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

$ws.Cells.Item(8, 2).Value = "XRP"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue 8 4 "0.596"
$ws.Cells.Item(8, 5).Value = "  -1.89%  "

Set-TextValue 9 4 "3.415.34"
$ws.Cells.Item(9, 5).Value = "  -2.59%  "

Set-TextValue 10 4 "7.16"
$ws.Cells.Item(10, 5).Value = "  -2.29%  "

Set-TextValue 11 4 "0.121"
$ws.Cells.Item(11, 5).Value = "  -3.79%  "

Set-TextValue 12 4 "0.436"
$ws.Cells.Item(12, 5).Value = "  -2.85%  "

Set-TextValue 13 4 "4.003.75"
$ws.Cells.Item(13, 5).Value = "  -2.50%  "

Set-TextValue 15 4 "0.0000186"
$ws.Cells.Item(15, 5).Value = "  -5.04%  "

Set-TextValue 16 4 "27.48"
$ws.Cells.Item(16, 5).Value = "  -4.74%  "

Set-TextValue 17 4 "64.555.83"
$ws.Cells.Item(17, 5).Value = "  -1.93%  "

Set-TextValue 18 4 "3.418.25"
$ws.Cells.Item(18, 5).Value = "  -2.37%  "

Set-TextValue 19 4 "6.31"
$ws.Cells.Item(19, 5).Value = "  -2.53%  "

Set-TextValue 20 4 "13.73"
$ws.Cells.Item(20, 5).Value = "  -4.29%  "

Set-TextValue 21 4 "377.86"
$ws.Cells.Item(21, 5).Value = "  -3.63%  "

Set-TextValue 22 4 "7.93"
$ws.Cells.Item(22, 5).Value = "  -4.60%  "

Set-TextValue 23 4 "0.544"
$ws.Cells.Item(23, 5).Value = "  -1.90%  "

Set-TextValue 24 4 "1.00"
$ws.Cells.Item(24, 5).Value = "  -0.44%  "

Set-TextValue 25 4 "71.77"
$ws.Cells.Item(25, 5).Value = "  -2.67%  "

Set-TextValue 26 4 "0.0000117"
$ws.Cells.Item(26, 5).Value = "  -6.39%  "

Set-TextValue 27 4 "9.87"
$ws.Cells.Item(27, 5).Value = "  +0.90%  "

$ws.Cells.Item(28, 5).Value = "  -1.30%  "

$ws.Cells.Item(29, 5).Value = "  +0.04%  "

Set-TextValue 30 4 "1.45"
$ws.Cells.Item(30, 5).Value = "  -0.59%  "

Set-TextValue 31 4 "6.14"
$ws.Cells.Item(31, 5).Value = "  -3.90%  "

$ws.Cells.Item(32, 5).Value = "  -3.28%  "

Set-TextValue 33 4 "23.11"
$ws.Cells.Item(33, 5).Value = "  -2.91%  "

Set-TextValue 34 4 "6.98"
$ws.Cells.Item(34, 5).Value = "  -2.91%  "

Set-TextValue 35 4 "1.56"
$ws.Cells.Item(35, 5).Value = "  +0.67%  "

Set-TextValue 36 4 "160.65"
$ws.Cells.Item(36, 5).Value = "  -1.46%  "

Set-TextValue 37 4 "1.89"
$ws.Cells.Item(37, 5).Value = "  -3.45%  "

Set-TextValue 38 4 "0.0750"
$ws.Cells.Item(38, 5).Value = "  -3.42%  "

Set-TextValue 39 4 "2.874.22"
$ws.Cells.Item(39, 5).Value = "  -7.57%  "

Set-TextValue 40 4 "6.67"
$ws.Cells.Item(40, 5).Value = "  +2.04%  "

Set-TextValue 41 4 "26.12"
$ws.Cells.Item(41, 5).Value = "  -4.87%  "

$ws.Cells.Item(42, 2).Value = "OKB"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 42 4 "43.02"
$ws.Cells.Item(42, 5).Value = "  -0.37%  "

$ws.Cells.Item(43, 2).Value = "Filecoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 43 4 "4.51"
$ws.Cells.Item(43, 5).Value = "  -0.91%  "

Set-TextValue 44 4 "0.0314"
$ws.Cells.Item(44, 5).Value = "  -3.11%  "

$ws.Cells.Item(45, 5).Value = "  -2.01%  "

Set-TextValue 46 4 "25.79"
$ws.Cells.Item(46, 5).Value = "  -0.56%  "

Set-TextValue 47 4 "313.57"
$ws.Cells.Item(47, 5).Value = "  -0.55%  "

$ws.Cells.Item(48, 2).Value = "ONDO"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue 48 4 "1.06"
$ws.Cells.Item(48, 5).Value = "  -5.53%  "

$ws.Cells.Item(49, 2).Value = "dogwifhat"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 49 4 "2.20"
$ws.Cells.Item(49, 5).Value = "  -2.65%  "

Set-TextValue 50 4 "0.106"
$ws.Cells.Item(50, 5).Value = "  -3.16%  "

Set-TextValue 51 4 "6.48"
$ws.Cells.Item(51, 5).Value = "  -3.78%  "

